$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The part described in row 9 changed from steel sheet stock to aluminum plate stock.
$ws.Range("A9").Value = "3/16 inch. thick`n6061 Aluminum Plate (2 ftX3 ft)"

# Its unit price changed; the shared formula in D9 and the SUM total in F2
# both recalculate automatically from this.
$ws.Range("C9").Value = 109.47

# The new two-line description wraps onto fewer lines than the old three-line
# one, so the row shrinks back down to the "two wrapped lines" height used by
# the similar row above it.
$ws.Rows.Item(9).RowHeight = 29.25

# Reflect the scrolled viewport / current selection recorded for the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("F2").Select()
